# Update need_to_buy.xlsx data: shift the date series forward by one day
# and refresh the forecast/quantity columns with the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(45904, 5629.64507100509, 5062.83913773471, 5688, 6232.053417, -0.948021511266234)
    3  = @(45905, 5629.6450710067, 4350.90531144263, 2952, 6232.05343, 83.3880696014971)
    4  = @(45906, 1140.44365045223, 2187.98669591926, 2952, 2189.567582, 11.8796094777929)
    5  = @(45907, 1020.77259022547, 2142.47378242945, 2952, 2061.465954, 9.63196442516583)
    6  = @(45908, 5912.85235826638, 4931.09520724043, 2952, 6520.373042, 107.775662123919)
    7  = @(45909, 6093.87809712743, 5162.81115716606, 2952, 6711.693705, 117.859448543276)
    8  = @(45910, 6093.87809712743, 5183.70297560104, 2952, 6711.693705, 118.729940978067)
    9  = @(45911, 6093.87809712743, 5170.31495142151, 2952, 6711.693705, 118.172106637254)
    10 = @(45912, 6093.87809712743, 4502.87079174772, 2952, 6711.693705, 90.3619333175121)
    11 = @(45913, 1155.32029950294, 2213.49521488796, 2952, 2179.639649, 11.9089401827092)
    12 = @(45914, 1023.06998026009, 2170.70210505867, 2952, 2036.153087, 9.65771715827403)
    13 = @(45915, 6061.80298841384, 5007.90996080768, 2952, 6722.522251, 113.192884308077)
    14 = @(45916, 6061.80298841384, 5179.63595066378, 2952, 6722.522251, 120.348133885414)
    15 = @(45917, 6061.80298841384, 5174.67833123251, 2952, 6722.522251, 120.141566409111)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
